$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '27.483.06'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.99%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.618.25'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.78%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '211.25'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.01%  '

# Row 6
$ws.Range("E6").Value = '  -1.33%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '22.85'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.78%  '

# Row 9
$ws.Range("E9").Value = '  +0.35%  '

# Row 10
$ws.Range("E10").Value = '  -0.31%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0888'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.42%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.846.56'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.78%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.620.72'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.56%  '

# Row 14
$ws.Range("E14").Value = '  -0.25%  '

# Row 15
$ws.Range("E15").Value = '  -2.45%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '27.461.88'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.92%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '233.18'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.55%  '

# Row 19
$ws.Range("E19").Value = '  -0.95%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.56'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.69%  '

# Row 21
$ws.Range("E21").Value = '  +0.07%  '

# Row 22
$ws.Range("E22").Value = '  -0.70%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.15'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.54%  '

# Row 24
$ws.Range("E24").Value = '  +5.89%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '150.82'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.58%  '

# Row 26
$ws.Range("E26").Value = '  -1.68%  '

# Row 27
$ws.Range("E27").Value = '  -1.17%  '

# Row 28
$ws.Range("E28").Value = '  +0.04%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.57'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '

# Row 30
$ws.Range("E30").Value = '  -0.83%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.0484'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.09%  '

# Row 32
$ws.Range("E32").Value = '  -1.28%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.472.35'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.76%  '

# Row 34
$ws.Range("E34").Value = '  -2.93%  '

# Row 35
$ws.Range("E35").Value = '  -3.18%  '

# Row 36
$ws.Range("E36").Value = '  -0.39%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.958'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +6.91%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0167'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.58%  '

# Row 39
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.558'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.39%  '

# Row 40
$ws.Range("E40").Value = '  -2.99%  '

# Row 41
$ws.Range("E41").Value = '  +0.06%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '68.03'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +2.47%  '

# Row 43
$ws.Range("B43").Value = 'mCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.46'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.40%  '

# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.983'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -4.80%  '

# Row 45
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.32%  '

# Row 46
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '5.27'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -7.67%  '

# Row 47
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.757.85'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.73%  '

# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.73'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.94%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '86.60'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.01%  '

# Row 50
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0₆0105'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.59%  '

# Row 51
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.101'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.61%  '
